# New data file for Nemegt (and related recompute of Noyon/Tost rows) for TNN
# Updates the density-estimate summary columns (L=D, M=SE, N=LCL, O=UCL) for the
# "model fit" rows and their comparison rows on Sheet1. Dependent formula
# columns (P/Q/R/S/T) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Noyon : row 2 (combined model) ---
$ws.Range("L2").Value = 14.2
$ws.Range("M2").Value = 0.46
$ws.Range("N2").Value = 14.01
$ws.Range("O2").Value = 16.88

# --- Noyon : row 6 (comparison model) ---
$ws.Range("L6").Value = 15.85
$ws.Range("M6").Value = 1.7
$ws.Range("N6").Value = 14.39
$ws.Range("O6").Value = 22.59

# --- Nemegt : row 14 (combined model) ---
$ws.Range("L14").Value = 20

# --- Nemegt : row 16 (comparison model) ---
$ws.Range("L16").Value = 18.510000000000002

# --- Tost : row 26 (combined model) ---
$ws.Range("L26").Value = 14.584
$ws.Range("M26").Value = 0.83
$ws.Range("N26").Value = 14.07
$ws.Range("O26").Value = 18.59

# --- Tost : row 31 (comparison model) ---
$ws.Range("L31").Value = 15.69
$ws.Range("M31").Value = 1.62
$ws.Range("N31").Value = 14.34
$ws.Range("O31").Value = 22.22

# --- Window / pane layout: split view with vertical split after column A,
# horizontal split after row 2, active pane bottom-right, with the final
# selection matching the refreshed view ---
$ws.Range("C1").Select() | Out-Null
$excel.ActiveWindow.SplitColumn = 1
$excel.ActiveWindow.SplitRow = 2
$excel.ActiveWindow.Split = $true

$ws.Range("A20").Select() | Out-Null
$ws.Range("A14").Select() | Out-Null

$excel.CalculateFull() | Out-Null
